$d = $word.ActiveDocument

# --- Rename the TOC bookmark on the "puzzles-cloud" Heading1 paragraph ---
$bm1 = $d.Bookmarks.Item("_Toc16364603065534479712892223")
$bm1Start = $bm1.Start
$bm1End = $bm1.End
$bm1.Delete()
$bm1Range = $d.Range($bm1Start, $bm1End)
$d.Bookmarks.Add("_Toc16364608742041690072513747", $bm1Range)

# --- Rename the TOC bookmark on the "mnogo je dobro bilo" Heading2 paragraph ---
$bm2 = $d.Bookmarks.Item("_Toc16364603065854199637577553")
$bm2Start = $bm2.Start
$bm2End = $bm2.End
$bm2.Delete()
$bm2Range = $d.Range($bm2Start, $bm2End)
$d.Bookmarks.Add("_Toc1636460874228703875388626", $bm2Range)

# --- Remove the entire "asd" Heading3 paragraph (bookmark, runs, marks and all) ---
$bm3 = $d.Bookmarks.Item("_Toc16364603066388608835155560")
$asdParagraph = $bm3.Range.Paragraphs(1)
$asdParagraph.Range.Delete()
